# Update "Riders" (col C) and "Average" (col D) figures on the Ridership
# sheet with the new Madigan bike hours data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Row 2
$ws.Range("C2").Value = 244
$ws.Range("D2").Value = 227.68

# Row 3
$ws.Range("C3").Value = 269
$ws.Range("D3").Value = 220.11

# Row 4
$ws.Range("C4").Value = 269
$ws.Range("D4").Value = 213.93

# Row 5
$ws.Range("C5").Value = 200

# Row 6
$ws.Range("C6").Value = 232
$ws.Range("D6").Value = 241.41

# Row 7
$ws.Range("C7").Value = 84
$ws.Range("D7").Value = 112.11

# Row 8
$ws.Range("C8").Value = 58
$ws.Range("D8").Value = 90.65000000000001

$wb.Save()
